$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.939.51'
$ws.Range("E2").Value = '  -1.15%  '

# Row 3
$ws.Range("D3").Value = '1.811.18'
$ws.Range("E3").Value = '  -0.61%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.64'
$ws.Range("E5").Value = '  -0.82%  '

# Row 6
$ws.Range("E6").Value = '  +0.08%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4615'
$ws.Range("E7").Value = '  +3.00%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3712'
$ws.Range("E8").Value = '  -1.27%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07374'
$ws.Range("E9").Value = '  -0.37%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8741'
$ws.Range("E10").Value = '  -0.55%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.48'
$ws.Range("E11").Value = '  -1.74%  '

# Row 12
$ws.Range("D12").Value = '1.795.02'
$ws.Range("E12").Value = '  -1.41%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.360'
$ws.Range("E13").Value = '  -1.02%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.36'
$ws.Range("E14").Value = '  -0.51%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.526'
$ws.Range("E15").Value = '  -2.78%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07049'
$ws.Range("E16").Value = '  -0.14%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  +0.17%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008700'
$ws.Range("E18").Value = '  -1.18%  '

# Row 19
$ws.Range("E19").Value = '  +0.05%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.72'
$ws.Range("E20").Value = '  -1.96%  '

# Row 21
$ws.Range("D21").Value = '26.947.15'
$ws.Range("E21").Value = '  -1.11%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.315'
$ws.Range("E22").Value = '  -0.48%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.64'
$ws.Range("E23").Value = '  -2.75%  '

# Row 24
$ws.Range("D24").Value = '2.013.31'
$ws.Range("E24").Value = '  -1.54%  '

# Row 25
$ws.Range("E25").Value = '  -3.04%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.40'
$ws.Range("E26").Value = '  +0.27%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.38'
$ws.Range("E27").Value = '  -1.04%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.141'
$ws.Range("E28").Value = '  -6.11%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.293'
$ws.Range("E29").Value = '  -0.89%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.92'
$ws.Range("E30").Value = '  -1.21%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08923'
$ws.Range("E31").Value = '  +0.42%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7520'
$ws.Range("E32").Value = '  -4.69%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.157'
$ws.Range("E33").Value = '  -3.28%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.924'
$ws.Range("E34").Value = '  -0.06%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.447'
$ws.Range("E35").Value = '  -2.74%  '

# Row 36
$ws.Range("E36").Value = '  +0.09%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.105'
$ws.Range("E37").Value = '  -0.20%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01977'
$ws.Range("E38").Value = '  +0.17%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05241'
$ws.Range("E39").Value = '  -0.33%  '

# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.428'
$ws.Range("E40").Value = '  +4.39%  '

# Row 41
$ws.Range("E41").Value = '  +1.95%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5311'
$ws.Range("E42").Value = '  +0.46%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.189'
$ws.Range("E43").Value = '  -1.32%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1665'
$ws.Range("E44").Value = '  -2.01%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.515'
$ws.Range("E45").Value = '  -1.33%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4984'
$ws.Range("E46").Value = '  -1.10%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.39'
$ws.Range("E47").Value = '  -2.02%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '104.23'
$ws.Range("E48").Value = '  -0.94%  '

# Row 49
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.001'
$ws.Range("E49").Value = '  +0.05%  '

# Row 50
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.669'
$ws.Range("E50").Value = '  -0.94%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06299'
$ws.Range("E51").Value = '  -1.32%  '
